# Update phase 2 tasks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase 2 Tasks")

# Reassign "I will take this" column from "Vivian" to "Team" for rows 7-10
$ws.Range("E7:E10").Value = "Team"

# Update the saved selection / active cell for this sheet's view
$ws.Activate()
$ws.Range("E22").Select()
